# Rename the original sheet to Sheet1 and add Sheet2 / Sheet3 after it,
# matching the target workbook's sheet list.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

$ws1.Activate()

# Overwrite the sample rows with the new learner data (headers in row 1
# are left as-is - same text, only their bold/fill styling changes below).
$ws1.Range("A2").Value = "Akshay"
$ws1.Range("B2").Value = "IT"
$ws1.Range("C2").Value = 7656787890
$ws1.Range("D2").Value = "akshay@gmail.com"

$ws1.Range("A3").Value = "Vinod"
$ws1.Range("B3").Value = "HR"
$ws1.Range("C3").Value = 7869352434
$ws1.Range("D3").Value = "vinod@gmail.com"

$ws1.Range("E3").Select()
